$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 18753544
$ws.Range("I40").Value = 5097.5
$ws.Range("K40").Value = 5097.5
$ws.Range("M40").Value = -4922.5
$ws.Range("H43").Value = 4283.9
$ws.Range("I43").Value = 2427.1428
$ws.Range("J43").Value = 8616.333000000001
$ws.Range("K43").Value = 2427.1428
$ws.Range("L43").Value = 8616.333000000001
$ws.Range("M43").Value = -2358.1428
$ws.Range("N43").Value = -8754.333000000001
$ws.Range("H62").Value = 7582734.5
$ws.Range("I62").Value = 11117829
$ws.Range("K62").Value = 11117829
$ws.Range("M62").Value = -11117205
$ws.Range("H65").Value = 7582734.5
$ws.Range("I65").Value = 11117829
$ws.Range("K65").Value = 55589145
$ws.Range("M65").Value = -55586025
$ws.Range("H88").Value = 2472.7144
$ws.Range("J88").Value = 2523.889
$ws.Range("L88").Value = 2523.889
$ws.Range("N88").Value = -3335.889
$ws.Range("H91").Value = 2472.7144
$ws.Range("J91").Value = 2523.889
$ws.Range("L91").Value = 2523.889
$ws.Range("N91").Value = -5331.889
$ws.Range("H107").Value = 2522.2
$ws.Range("I107").Value = 2462.1667
$ws.Range("J107").Value = 2612.25
$ws.Range("K107").Value = 2462.1667
$ws.Range("L107").Value = 2612.25
$ws.Range("M107").Value = -542.1667000000002
$ws.Range("N107").Value = -6452.25
$ws.Range("H113").Value = 4553.5557
$ws.Range("J113").Value = 5071.143
$ws.Range("L113").Value = 5071.143
$ws.Range("N113").Value = -11579.143
$ws.Range("H132").Value = 3691.1924
$ws.Range("J132").Value = 3746.25
$ws.Range("L132").Value = 11238.75
$ws.Range("N132").Value = -16298.75
$ws.Range("H138").Value = 3521.6428
$ws.Range("I138").Value = 1881
$ws.Range("J138").Value = 3795.0833
$ws.Range("K138").Value = 5643
$ws.Range("L138").Value = 11385.2499
$ws.Range("M138").Value = -503
$ws.Range("N138").Value = -21665.2499

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1905.6923
$ws.Range("I2").Value = 1964.9166
$ws.Range("K2").Value = 1964.9166
$ws.Range("M2").Value = -1851.9166
$ws.Range("H32").Value = 12139.854
$ws.Range("I32").Value = 7748
$ws.Range("K32").Value = 7748
$ws.Range("M32").Value = -7461
$ws.Range("H61").Value = 6862.8887
$ws.Range("H63").Value = 2903.4
$ws.Range("I63").Value = 2580.5293
$ws.Range("K63").Value = 2580.5293
$ws.Range("M63").Value = -1894.5293
$ws.Range("H66").Value = 2903.4
$ws.Range("I66").Value = 2580.5293
$ws.Range("K66").Value = 12902.6465
$ws.Range("M66").Value = -9470.646500000001
$ws.Range("H102").Value = 2969.8333
$ws.Range("I102").Value = 2267.5
$ws.Range("J102").Value = 4374.5
$ws.Range("K102").Value = 2267.5
$ws.Range("L102").Value = 4374.5
$ws.Range("M102").Value = -645.5
$ws.Range("N102").Value = -7618.5
$ws.Range("H116").Value = 1905.6923
$ws.Range("I116").Value = 1964.9166
$ws.Range("K116").Value = 1964.9166
$ws.Range("M116").Value = 329.0834
$ws.Range("H122").Value = 2233.4167
$ws.Range("I122").Value = 2362.5
$ws.Range("J122").Value = 1846.1666
$ws.Range("K122").Value = 7087.5
$ws.Range("L122").Value = 5538.4998
$ws.Range("M122").Value = -4637.5
$ws.Range("N122").Value = -10438.4998
$ws.Range("H136").Value = 6862.8887

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1905.6923
$ws.Range("I3").Value = 1964.9166
$ws.Range("K3").Value = 1964.9166
$ws.Range("M3").Value = -1850.9166
$ws.Range("H20").Value = 1737.5714
$ws.Range("I20").Value = 1324.2858
$ws.Range("K20").Value = 1324.2858
$ws.Range("M20").Value = -1077.2858
$ws.Range("H86").Value = 3670.182
$ws.Range("I86").Value = 3110.2856
$ws.Range("K86").Value = 3110.2856
$ws.Range("M86").Value = -1987.2856
$ws.Range("H89").Value = 3670.182
$ws.Range("I89").Value = 3110.2856
$ws.Range("K89").Value = 15551.428
$ws.Range("M89").Value = -9935.428
$ws.Range("H107").Value = 3164.0588
$ws.Range("I107").Value = 3185.7856
$ws.Range("K107").Value = 3185.7856
$ws.Range("M107").Value = -1265.7856
$ws.Range("H134").Value = 8376.25
$ws.Range("I134").Value = 8199.4
$ws.Range("J134").Value = 8671
$ws.Range("K134").Value = 24598.2
$ws.Range("L134").Value = 26013
$ws.Range("M134").Value = -22063.2
$ws.Range("N134").Value = -31083

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2971.5833
$ws.Range("I122").Value = 2969.0908
$ws.Range("J122").Value = 2999
$ws.Range("K122").Value = 8907.2724
$ws.Range("L122").Value = 8997
$ws.Range("M122").Value = -6457.2724
$ws.Range("N122").Value = -13897
$ws.Range("H134").Value = 2938.875
$ws.Range("I134").Value = 2892.739
$ws.Range("K134").Value = 8678.217000000001
$ws.Range("M134").Value = -6143.217000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 116.55556
$ws.Range("I2").Value = 83
$ws.Range("J2").Value = 143.4
$ws.Range("K2").Value = 498
$ws.Range("L2").Value = 860.4000000000001
$ws.Range("M2").Value = -385
$ws.Range("N2").Value = -1086.4
$ws.Range("H23").Value = 897.375
$ws.Range("I23").Value = 137.4
$ws.Range("J23").Value = 1242.8182
$ws.Range("K23").Value = 412.2
$ws.Range("L23").Value = 3728.4546
$ws.Range("M23").Value = -177.2
$ws.Range("N23").Value = -4198.4546

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 130.5
$ws.Range("I2").Value = 117.85714
$ws.Range("J2").Value = 160
$ws.Range("K2").Value = 117.85714
$ws.Range("L2").Value = 160
$ws.Range("M2").Value = -4.857140000000001
$ws.Range("N2").Value = -386
$ws.Range("H80").Value = 3995.639
$ws.Range("I80").Value = 3739.6316
$ws.Range("J80").Value = 4281.7646
$ws.Range("K80").Value = 3739.6316
$ws.Range("L80").Value = 4281.7646
$ws.Range("M80").Value = -2741.6316
$ws.Range("N80").Value = -6277.7646
$ws.Range("H83").Value = 3995.639
$ws.Range("I83").Value = 3739.6316
$ws.Range("J83").Value = 4281.7646
$ws.Range("K83").Value = 18698.158
$ws.Range("L83").Value = 21408.823
$ws.Range("M83").Value = -13706.158
$ws.Range("N83").Value = -31392.823
$ws.Range("H122").Value = 6416.8716
$ws.Range("I122").Value = 4805.263
$ws.Range("J122").Value = 7947.9
$ws.Range("K122").Value = 14415.789
$ws.Range("L122").Value = 23843.7
$ws.Range("M122").Value = -11965.789
$ws.Range("N122").Value = -28743.7

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1688.7609
$ws.Range("I16").Value = 1662.5588
$ws.Range("J16").Value = 1763
$ws.Range("K16").Value = 1662.5588
$ws.Range("L16").Value = 1763
$ws.Range("M16").Value = -1492.5588
$ws.Range("N16").Value = -2103
$ws.Range("H31").Value = 785.5714
$ws.Range("I31").Value = 785.5714
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 785.5714
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -537.5714
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("H61").Value = 2824.45
$ws.Range("I61").Value = 2891.2307
$ws.Range("J61").Value = 2700.4285
$ws.Range("K61").Value = 2891.2307
$ws.Range("L61").Value = 2700.4285
$ws.Range("M61").Value = -2689.2307
$ws.Range("N61").Value = -3104.4285
$ws.Range("H68").Value = 8799.75
$ws.Range("I68").Value = 9157.632
$ws.Range("K68").Value = 9157.632
$ws.Range("M68").Value = -8408.632
$ws.Range("H71").Value = 8799.75
$ws.Range("I71").Value = 9157.632
$ws.Range("K71").Value = 45788.16
$ws.Range("M71").Value = -42044.16
$ws.Range("H82").Value = 3111.111
$ws.Range("I82").Value = 1922.8462
$ws.Range("K82").Value = 1922.8462
$ws.Range("M82").Value = -1561.8462
$ws.Range("H85").Value = 3111.111
$ws.Range("I85").Value = 1922.8462
$ws.Range("K85").Value = 1922.8462
$ws.Range("M85").Value = -674.8462
$ws.Range("H113").Value = 2824.45
$ws.Range("I113").Value = 2891.2307
$ws.Range("J113").Value = 2700.4285
$ws.Range("K113").Value = 2891.2307
$ws.Range("L113").Value = 2700.4285
$ws.Range("M113").Value = -721.2307000000001
$ws.Range("N113").Value = -7040.4285
$ws.Range("H122").Value = 4226.1
$ws.Range("I122").Value = 4130.125
$ws.Range("K122").Value = 12390.375
$ws.Range("M122").Value = -9940.375
$ws.Range("H132").Value = 3641
$ws.Range("I132").Value = 3703.158
$ws.Range("K132").Value = 11109.474
$ws.Range("M132").Value = -8579.474
$ws.Range("N31").ClearContents()
$ws.Range("M41").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 4395.4443
$ws.Range("I4").Value = 3344.3333
$ws.Range("K4").Value = 3344.3333
$ws.Range("M4").Value = -3231.3333
$ws.Range("H81").Value = 2290.4707
$ws.Range("H84").Value = 2290.4707
$ws.Range("H96").Value = 37291.07
$ws.Range("J96").Value = 2754.1667
$ws.Range("L96").Value = 2754.1667
$ws.Range("N96").Value = -5500.1667
$ws.Range("H122").Value = 2457.2856
$ws.Range("I122").Value = 2461.5
$ws.Range("J122").Value = 2432
$ws.Range("K122").Value = 7384.5
$ws.Range("L122").Value = 7296
$ws.Range("M122").Value = -4934.5
$ws.Range("N122").Value = -12196
